$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    if ($value -match '^-?\d+(\.\d+)?$') {
        $r.Value = "'" + $value
        $r.Style = "Normal"
    } else {
        $r.Value = $value
    }
}

Set-TextValue "D2" "60.242.01"
Set-TextValue "E2" "  +2.82%  "
Set-TextValue "D3" "2.319.89"
Set-TextValue "E3" "  +0.82%  "
Set-TextValue "E4" "  +0.14%  "
Set-TextValue "D5" "545.06"
Set-TextValue "E5" "  +1.07%  "
Set-TextValue "D6" "131.02"
Set-TextValue "E6" "  -1.25%  "
Set-TextValue "E7" "  +0.12%  "
Set-TextValue "D8" "0.580"
Set-TextValue "E8" "  -1.49%  "
Set-TextValue "D9" "2.316.93"
Set-TextValue "E9" "  +0.99%  "
Set-TextValue "E10" "  -0.31%  "
Set-TextValue "D11" "5.51"
Set-TextValue "E12" "  -0.07%  "
Set-TextValue "D13" "0.334"
Set-TextValue "E13" "  -0.35%  "
Set-TextValue "D14" "23.68"
Set-TextValue "E14" "  -1.09%  "
Set-TextValue "D15" "60.200.37"
Set-TextValue "E15" "  +2.91%  "
Set-TextValue "D16" "2.732.85"
Set-TextValue "E16" "  +1.04%  "
Set-TextValue "E17" "  -0.47%  "
Set-TextValue "D18" "2.310.32"
Set-TextValue "E18" "  +1.04%  "
Set-TextValue "E19" "  -0.44%  "
Set-TextValue "E20" "  -2.01%  "
Set-TextValue "D21" "313.96"
Set-TextValue "E21" "  -0.95%  "
Set-TextValue "D22" "6.64"
Set-TextValue "E22" "  +0.67%  "
Set-TextValue "D23" "0.997"
Set-TextValue "E23" "  -0.34%  "
Set-TextValue "D24" "63.66"
Set-TextValue "E24" "  +0.50%  "
Set-TextValue "E25" "  +1.76%  "
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  +0.10%  "
Set-TextValue "E27" "  -2.41%  "
Set-TextValue "E28" "  +2.55%  "
Set-TextValue "D29" "173.41"
Set-TextValue "E29" "  +1.31%  "
Set-TextValue "E30" "  +8.18%  "
Set-TextValue "E31" "  +1.08%  "
Set-TextValue "E32" "  -0.37%  "
Set-TextValue "D33" "5.93"
Set-TextValue "E33" "  +1.07%  "
Set-TextValue "D34" "1.36"
Set-TextValue "E34" "  +9.01%  "
Set-TextValue "E35" "  -1.20%  "
Set-TextValue "E36" "  +0.00%  "
Set-TextValue "E37" "  -0.73%  "
Set-TextValue "E38" "  +0.09%  "
Set-TextValue "E39" "  +1.51%  "
Set-TextValue "D40" "323.23"
Set-TextValue "E40" "  +9.94%  "
Set-TextValue "E41" "  -1.13%  "
Set-TextValue "E42" "  +0.84%  "
Set-TextValue "D43" "138.22"
Set-TextValue "E43" "  -2.16%  "
Set-TextValue "D44" "3.48"
Set-TextValue "E44" "  +0.50%  "
Set-TextValue "E45" "  -1.50%  "
Set-TextValue "D46" "19.15"
Set-TextValue "E46" "  +4.06%  "
Set-TextValue "D47" "0.0495"
Set-TextValue "E47" "  -0.59%  "
Set-TextValue "D48" "0.559"
Set-TextValue "E48" "  +0.28%  "
Set-TextValue "E49" "  +0.29%  "
Set-TextValue "E50" "  +16.50%  "
Set-TextValue "D51" "11.02"

Write-Host "Applied cryptos list update"
